$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-27 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-28 Friday", 2) | Out-Null
$d.Content.Find.Execute("175÷9=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "288÷9=32, 0", 2) | Out-Null
$d.Content.Find.Execute("480÷3=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "898÷8=112, 2", 2) | Out-Null
$d.Content.Find.Execute("272÷6=45, 2", $true, $false, $false, $false, $false, $true, 1, $false, "424÷4=106, 0", 2) | Out-Null
$d.Content.Find.Execute("531÷5=106, 1", $true, $false, $false, $false, $false, $true, 1, $false, "130÷7=18, 4", 2) | Out-Null
$d.Content.Find.Execute("928÷2=464, 0", $true, $false, $false, $false, $false, $true, 1, $false, "230÷7=32, 6", 2) | Out-Null
$d.Content.Find.Execute("782÷8=97, 6", $true, $false, $false, $false, $false, $true, 1, $false, "355÷2=177, 1", 2) | Out-Null
$d.Content.Find.Execute("663÷4=165, 3", $true, $false, $false, $false, $false, $true, 1, $false, "361÷8=45, 1", 2) | Out-Null
$d.Content.Find.Execute("417÷5=83, 2", $true, $false, $false, $false, $false, $true, 1, $false, "269÷4=67, 1", 2) | Out-Null
$d.Content.Find.Execute("588÷7=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "552÷4=138, 0", 2) | Out-Null
$d.Content.Find.Execute("814÷8=101, 6", $true, $false, $false, $false, $false, $true, 1, $false, "584÷8=73, 0", 2) | Out-Null
$d.Content.Find.Execute("797÷9=88, 5", $true, $false, $false, $false, $false, $true, 1, $false, "310÷8=38, 6", 2) | Out-Null
$d.Content.Find.Execute("545÷8=68, 1", $true, $false, $false, $false, $false, $true, 1, $false, "530÷3=176, 2", 2) | Out-Null
$d.Content.Find.Execute("147÷4=36, 3", $true, $false, $false, $false, $false, $true, 1, $false, "674÷6=112, 2", 2) | Out-Null
$d.Content.Find.Execute("227÷8=28, 3", $true, $false, $false, $false, $false, $true, 1, $false, "245÷3=81, 2", 2) | Out-Null
$d.Content.Find.Execute("891÷7=127, 2", $true, $false, $false, $false, $false, $true, 1, $false, "699÷4=174, 3", 2) | Out-Null
$d.Content.Find.Execute("764÷2=382, 0", $true, $false, $false, $false, $false, $true, 1, $false, "640÷3=213, 1", 2) | Out-Null
$d.Content.Find.Execute("650÷7=92, 6", $true, $false, $false, $false, $false, $true, 1, $false, "926÷2=463, 0", 2) | Out-Null
$d.Content.Find.Execute("903÷3=301, 0", $true, $false, $false, $false, $false, $true, 1, $false, "363÷9=40, 3", 2) | Out-Null
$d.Content.Find.Execute("307÷6=51, 1", $true, $false, $false, $false, $false, $true, 1, $false, "286÷4=71, 2", 2) | Out-Null
$d.Content.Find.Execute("680÷7=97, 1", $true, $false, $false, $false, $false, $true, 1, $false, "285÷6=47, 3", 2) | Out-Null
$d.Content.Find.Execute("641÷4=160, 1", $true, $false, $false, $false, $false, $true, 1, $false, "670÷9=74, 4", 2) | Out-Null
$d.Content.Find.Execute("706÷8=88, 2", $true, $false, $false, $false, $false, $true, 1, $false, "182÷2=91, 0", 2) | Out-Null
$d.Content.Find.Execute("212÷3=70, 2", $true, $false, $false, $false, $false, $true, 1, $false, "454÷2=227, 0", 2) | Out-Null
$d.Content.Find.Execute("425÷9=47, 2", $true, $false, $false, $false, $false, $true, 1, $false, "131÷7=18, 5", 2) | Out-Null
$d.Content.Find.Execute("249÷5=49, 4", $true, $false, $false, $false, $false, $true, 1, $false, "912÷5=182, 2", 2) | Out-Null
